$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma-separated person names: commas used as list separators
#     were mistakenly also matching decimal-style replacement; normalize
#     them to periods (and drop stray internal periods in abbreviations).
$ws.Range("E42").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E171").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E88").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E184").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Fix "Importe" (amount) column: values were scraped in
#     Argentine/Spanish locale formatting ("." thousands separator,
#     "," decimal separator). Re-write them in plain dot-decimal
#     notation without the thousands separator, keeping them as text.
$ws.Range("H2:H248").NumberFormat = "@"

$importe = @{
  2 = "750.00"
  3 = "4460.00"
  4 = "35000.00"
  5 = "689500.00"
  6 = "272000.00"
  7 = "35000.00"
  8 = "70000.00"
  9 = "79815.00"
  10 = "25830.00"
  11 = "744.54"
  12 = "850.00"
  13 = "4998.15"
  14 = "150441.82"
  15 = "28645.20"
  16 = "45673.46"
  17 = "6970.00"
  18 = "7290.00"
  19 = "20634.73"
  20 = "5224.00"
  21 = "27827.33"
  22 = "12608.00"
  23 = "3600.00"
  24 = "4500.00"
  25 = "9190.00"
  26 = "6890.00"
  27 = "31.36"
  28 = "2626.98"
  29 = "30.00"
  30 = "269239.92"
  31 = "14098.35"
  32 = "3916.00"
  33 = "5545.77"
  34 = "14803.00"
  35 = "3038.10"
  36 = "1769.00"
  37 = "746.00"
  38 = "17793.94"
  39 = "3110.74"
  40 = "407.80"
  41 = "39703.57"
  42 = "19124.00"
  43 = "7971.22"
  44 = "5863.22"
  45 = "386118.60"
  46 = "7001.69"
  47 = "1047.00"
  48 = "20870.85"
  49 = "221.48"
  50 = "775.87"
  51 = "343200.00"
  52 = "22449.89"
  53 = "185.00"
  54 = "7908.48"
  55 = "15110.00"
  56 = "250.00"
  57 = "1950.00"
  58 = "97182.00"
  59 = "12522.78"
  60 = "172.24"
  61 = "308.93"
  62 = "14431.14"
  63 = "600.00"
  64 = "4715.00"
  65 = "19740.00"
  66 = "366662.00"
  67 = "5243.23"
  68 = "1014.60"
  69 = "501.00"
  70 = "501.43"
  71 = "20940.00"
  72 = "460.00"
  73 = "21264.00"
  74 = "3337.80"
  75 = "150.00"
  76 = "294.00"
  77 = "1392.00"
  78 = "5169.00"
  79 = "7500.00"
  80 = "5660.00"
  81 = "57070.00"
  82 = "9067.00"
  83 = "117450.00"
  84 = "30560.00"
  85 = "104.50"
  86 = "7130.00"
  87 = "320.00"
  88 = "8522.00"
  89 = "4540.00"
  90 = "399.98"
  91 = "10400.00"
  92 = "577.80"
  93 = "2440.00"
  94 = "1500.00"
  95 = "526097.59"
  96 = "53261.49"
  97 = "8000.00"
  98 = "7.41"
  99 = "26.25"
  100 = "27994.92"
  101 = "5635.82"
  102 = "3650.00"
  103 = "1034.90"
  104 = "20760.00"
  105 = "1905.60"
  106 = "1950.00"
  107 = "1200.00"
  108 = "660.00"
  109 = "1497.00"
  110 = "1111.00"
  111 = "7810.00"
  112 = "22423.70"
  113 = "3328.00"
  114 = "7065.00"
  115 = "425.40"
  116 = "69.70"
  117 = "4800.00"
  118 = "8034.00"
  119 = "1253.00"
  120 = "1406.60"
  121 = "268.09"
  122 = "88600.00"
  123 = "10080.00"
  124 = "500.00"
  125 = "7000.00"
  126 = "1500.00"
  127 = "36000.00"
  128 = "2430.00"
  129 = "8800.00"
  130 = "36100.00"
  131 = "31460.00"
  132 = "9600.00"
  133 = "5500.00"
  134 = "621.25"
  135 = "427.50"
  136 = "58.00"
  137 = "195.65"
  138 = "16313.94"
  139 = "575.00"
  140 = "669.47"
  141 = "41.67"
  142 = "301500.00"
  143 = "8000.00"
  144 = "4900.00"
  145 = "10000.00"
  146 = "19035.72"
  147 = "2500.00"
  148 = "3204.50"
  149 = "2556.00"
  150 = "3000.00"
  151 = "3000.00"
  152 = "2000.00"
  153 = "1500.00"
  154 = "7500.00"
  155 = "4000.00"
  156 = "4100.00"
  157 = "17700.00"
  158 = "15000.00"
  159 = "6000.00"
  160 = "2500.00"
  161 = "1500.00"
  162 = "8640.00"
  163 = "1680.00"
  164 = "1300.00"
  165 = "3000.00"
  166 = "27063.80"
  167 = "95.20"
  168 = "4330.00"
  169 = "3000.00"
  170 = "480.04"
  171 = "40506.00"
  172 = "1208.50"
  173 = "1500.00"
  174 = "9297.00"
  175 = "48.24"
  176 = "196.65"
  177 = "2942.00"
  178 = "3500.00"
  179 = "2882.68"
  180 = "1776.79"
  181 = "71600.00"
  182 = "26936.65"
  183 = "9214.04"
  184 = "26380.00"
  185 = "520.00"
  186 = "5070.00"
  187 = "4280.00"
  188 = "2411.09"
  189 = "690.00"
  190 = "1665.62"
  191 = "50682.10"
  192 = "13307.25"
  193 = "42750.00"
  194 = "340.00"
  195 = "1200.00"
  196 = "10320.00"
  197 = "47.73"
  198 = "467220.13"
  199 = "25000.00"
  200 = "25000.00"
  201 = "8500.00"
  202 = "50000.00"
  203 = "25000.00"
  204 = "25000.00"
  205 = "50000.00"
  206 = "50000.00"
  207 = "25000.00"
  208 = "9200.00"
  209 = "2153596.16"
  210 = "10500.00"
  211 = "116190.00"
  212 = "122190.00"
  213 = "116190.00"
  214 = "116190.00"
  215 = "116190.00"
  216 = "116190.00"
  217 = "200190.00"
  218 = "200190.00"
  219 = "294690.00"
  220 = "116190.00"
  221 = "116190.00"
  222 = "116190.00"
  223 = "116190.00"
  224 = "116190.00"
  225 = "200190.00"
  226 = "369190.00"
  227 = "200190.00"
  228 = "116190.00"
  229 = "184190.00"
  230 = "116190.00"
  231 = "116190.00"
  232 = "116190.00"
  233 = "369961.74"
  234 = "11500333.74"
  235 = "15450.00"
  236 = "19801.54"
  237 = "181620.00"
  238 = "24786.00"
  239 = "3500.00"
  240 = "1238.42"
  241 = "3425.00"
  242 = "3500.00"
  243 = "16800.00"
  244 = "17000.00"
  245 = "105500.00"
  246 = "3500.00"
  247 = "8.80"
  248 = "1540.00"
}

foreach ($row in $importe.Keys) {
  $ws.Cells.Item($row, 8).Value = $importe[$row]
}
